$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 3 through 9 (old extra data rows), keeping header (row 1) and row 2
$ws.Range("A3:F9").EntireRow.Delete()

# Update row 2 with the new data
$ws.Range("A2").Value = "valdirene"
$ws.Range("B2").Value = "150,00"
$ws.Range("C2").Value = "0,00"
$ws.Range("D2").Value = ""
$ws.Range("E2").Value = "23/02/2026 01:15:46"
$ws.Range("F2").Value = "ADMINISTRADOR"
